# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the first (bc666980...) row
# on each sheet, reflecting a newly-regenerated handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for bc666980... row
$overview.Range("G2").Value = "2016-09-03 07:09:27"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
$zhcn.Range("H2").Value = "2016-09-03 07:09:22"
$zhcn.Range("K2").Value = "2016-09-03 07:09:40"

# de-de sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
$dede.Range("H2").Value = "2016-09-03 07:09:27"
$dede.Range("K2").Value = "2016-09-03 07:09:47"
